# feat: Adicao planilha disponibilidade dos professores
#
# The last two data rows of the schedule (row 29 "Trabalho de Conclusao de
# Curso II" / Matheus Fontanelle Pereira, and row 30 "Atividades de Extensao V"
# / Fernando da Silva Osorio) are removed from the table, turning row 28
# ("Atividades de Extensao III" / Natalia Madalena Boelter) into the new
# last row of the table (and picking up a closing border), while rows 29
# and 30 become blank spacer rows (content cleared, borders removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 becomes the final row of the table: B28/C28/D28 gain the same
# closing right border that A28 and E28 already have (thin, black,
# continuous). BorderAround is used (rather than poking the individual
# Borders.Item(...) properties) so the runtime recomputes each cell's
# style in a single shot instead of leaving unused transient styles
# behind in the workbook's style table.
$ws.Range("C28").BorderAround(1, 2, -4105, 0)   # xlContinuous, xlThin
$ws.Range("D28").BorderAround(1, 2, -4105, 0)
$ws.Range("B28").BorderAround(1, 2, -4105, 0)

# Clear the now-removed rows 29 and 30 (values + shared-string references).
# (Multi-area ranges only affect their first area when cleared/restyled in
# this runtime, so each row is handled individually.)
$row29 = $ws.Range("A29:E29")
$row30 = $ws.Range("A30:E30")

$row29.ClearContents()
$row30.ClearContents()

# Remove the borders that used to outline those rows, leaving plain cells.
$row29.Borders.LineStyle = -4142          # xlLineStyleNone
$row30.Borders.LineStyle = -4142          # xlLineStyleNone

Write-Host "Cleared rows 29-30 and closed the border on row 28"
